$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-13 in place with the recomputed TPM-based NATMI values
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tslp"
$ws.Range("C2").Value = "Il7r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.142057
$ws.Range("H2").Value = 3.426171
$ws.Range("I2").Value = 0.2487160836358648
$ws.Range("J2").Value = 0.2487160836358647
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.07349133333333334
$ws.Range("N2").Value = 0.220474
$ws.Range("O2").Value = 0.00285456188830886
$ws.Range("P2").Value = 0.00285456188830886
$ws.Range("Q2").Value = 0.08393129167266668
$ws.Range("R2").Value = 0.7553816250540001
$ws.Range("S2").Value = 0.0007099754533563784
$ws.Range("T2").Value = 0.0007099754533563783

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tslp"
$ws.Range("C3").Value = "Il7r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.142057
$ws.Range("H3").Value = 3.426171
$ws.Range("I3").Value = 0.2487160836358648
$ws.Range("J3").Value = 0.2487160836358647
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.178101
$ws.Range("N3").Value = 0.534303
$ws.Range("O3").Value = 0.006917826957414881
$ws.Range("P3").Value = 0.006917826957414882
$ws.Range("Q3").Value = 0.203401493757
$ws.Range("R3").Value = 1.830613443813
$ws.Range("S3").Value = 0.001720574828118839
$ws.Range("T3").Value = 0.001720574828118839

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tslp"
$ws.Range("C4").Value = "Il7r"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.142057
$ws.Range("H4").Value = 3.426171
$ws.Range("I4").Value = 0.2487160836358648
$ws.Range("J4").Value = 0.2487160836358647
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 25.493631
$ws.Range("N4").Value = 76.48089300000001
$ws.Range("O4").Value = 0.9902276111542762
$ws.Range("P4").Value = 0.9902276111542762
$ws.Range("Q4").Value = 29.11517973896701
$ws.Range("R4").Value = 262.036617650703
$ws.Range("S4").Value = 0.2462855333543895
$ws.Range("T4").Value = 0.2462855333543895

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tslp"
$ws.Range("C5").Value = "Il7r"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.159282
$ws.Range("H5").Value = 3.477846
$ws.Range("I5").Value = 0.2524673276986635
$ws.Range("J5").Value = 0.2524673276986635
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.07349133333333334
$ws.Range("N5").Value = 0.220474
$ws.Range("O5").Value = 0.00285456188830886
$ws.Range("P5").Value = 0.00285456188830886
$ws.Range("Q5").Value = 0.08519717988933334
$ws.Range("R5").Value = 0.766774619004
$ws.Range("S5").Value = 0.0007206836116917885
$ws.Range("T5").Value = 0.0007206836116917885

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tslp"
$ws.Range("C6").Value = "Il7r"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.159282
$ws.Range("H6").Value = 3.477846
$ws.Range("I6").Value = 0.2524673276986635
$ws.Range("J6").Value = 0.2524673276986635
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.178101
$ws.Range("N6").Value = 0.534303
$ws.Range("O6").Value = 0.006917826957414881
$ws.Range("P6").Value = 0.006917826957414882
$ws.Range("Q6").Value = 0.206469283482
$ws.Range("R6").Value = 1.858223551338
$ws.Range("S6").Value = 0.001746525285420311
$ws.Range("T6").Value = 0.001746525285420311

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tslp"
$ws.Range("C7").Value = "Il7r"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.159282
$ws.Range("H7").Value = 3.477846
$ws.Range("I7").Value = 0.2524673276986635
$ws.Range("J7").Value = 0.2524673276986635
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 25.493631
$ws.Range("N7").Value = 76.48089300000001
$ws.Range("O7").Value = 0.9902276111542762
$ws.Range("P7").Value = 0.9902276111542762
$ws.Range("Q7").Value = 29.554307532942
$ws.Range("R7").Value = 265.988767796478
$ws.Range("S7").Value = 0.2500001188015514
$ws.Range("T7").Value = 0.2500001188015514

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Tslp"
$ws.Range("C8").Value = "Il7r"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.894276
$ws.Range("H8").Value = 2.682828
$ws.Range("I8").Value = 0.1947545739044081
$ws.Range("J8").Value = 0.194754573904408
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.07349133333333334
$ws.Range("N8").Value = 0.220474
$ws.Range("O8").Value = 0.00285456188830886
$ws.Range("P8").Value = 0.00285456188830886
$ws.Range("Q8").Value = 0.06572153560800001
$ws.Range("R8").Value = 0.5914938204719999
$ws.Range("S8").Value = 0.0005559389842413545
$ws.Range("T8").Value = 0.0005559389842413544

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Tslp"
$ws.Range("C9").Value = "Il7r"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.894276
$ws.Range("H9").Value = 2.682828
$ws.Range("I9").Value = 0.1947545739044081
$ws.Range("J9").Value = 0.194754573904408
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.178101
$ws.Range("N9").Value = 0.534303
$ws.Range("O9").Value = 0.006917826957414881
$ws.Range("P9").Value = 0.006917826957414882
$ws.Range("Q9").Value = 0.159271449876
$ws.Range("R9").Value = 1.433443048884
$ws.Range("S9").Value = 0.001347278441435763
$ws.Range("T9").Value = 0.001347278441435763

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Tslp"
$ws.Range("C10").Value = "Il7r"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.894276
$ws.Range("H10").Value = 2.682828
$ws.Range("I10").Value = 0.1947545739044081
$ws.Range("J10").Value = 0.194754573904408
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 25.493631
$ws.Range("N10").Value = 76.48089300000001
$ws.Range("O10").Value = 0.9902276111542762
$ws.Range("P10").Value = 0.9902276111542762
$ws.Range("Q10").Value = 22.798342356156
$ws.Range("R10").Value = 205.185081205404
$ws.Range("S10").Value = 0.1928513564787309
$ws.Range("T10").Value = 0.1928513564787309

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Tslp"
$ws.Range("C11").Value = "Il7r"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.396195
$ws.Range("H11").Value = 4.188585
$ws.Range("I11").Value = 0.3040620147610637
$ws.Range("J11").Value = 0.3040620147610637
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.07349133333333334
$ws.Range("N11").Value = 0.220474
$ws.Range("O11").Value = 0.00285456188830886
$ws.Range("P11").Value = 0.00285456188830886
$ws.Range("Q11").Value = 0.1026082321433333
$ws.Range("R11").Value = 0.92347408929
$ws.Range("S11").Value = 0.0008679638390193384
$ws.Range("T11").Value = 0.0008679638390193384

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Tslp"
$ws.Range("C12").Value = "Il7r"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.396195
$ws.Range("H12").Value = 4.188585
$ws.Range("I12").Value = 0.3040620147610637
$ws.Range("J12").Value = 0.3040620147610637
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.178101
$ws.Range("N12").Value = 0.534303
$ws.Range("O12").Value = 0.006917826957414881
$ws.Range("P12").Value = 0.006917826957414882
$ws.Range("Q12").Value = 0.2486637256949999
$ws.Range("R12").Value = 2.237973531255
$ws.Range("S12").Value = 0.002103448402439968
$ws.Range("T12").Value = 0.002103448402439968

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Tslp"
$ws.Range("C13").Value = "Il7r"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.396195
$ws.Range("H13").Value = 4.188585
$ws.Range("I13").Value = 0.3040620147610637
$ws.Range("J13").Value = 0.3040620147610637
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 25.493631
$ws.Range("N13").Value = 76.48089300000001
$ws.Range("O13").Value = 0.9902276111542762
$ws.Range("P13").Value = 0.9902276111542762
$ws.Range("Q13").Value = 35.594080134045
$ws.Range("R13").Value = 320.346721206405
$ws.Range("S13").Value = 0.3010906025196043
$ws.Range("T13").Value = 0.3010906025196043

# Rows 14-17 no longer exist in the updated dataset (MuSCs target-cluster rows removed); delete them
$ws.Range("A14:T17").Delete()
